# Publication 0.2.0 preparation edit for
# CodeSystem-eclaire-type-contact-code-system.xlsx
#
# Changes:
#   1. "Metadata" sheet, Version value: 0.1.1 -> 0.2.0
#   2. "Metadata" sheet, Date value: 2023-10-19T16:17:18+00:00 -> 2023-10-19T17:05:12+00:00
#   3. "Metadata" sheet: new "Jurisdiction" / "iso:code:3166:FR" row inserted
#      right after the existing "Contact" row (pushes every following row down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update the Version property value.
$ws.Range("B3").Value = "0.2.0"

# 2. Update the Date property value.
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# 3. Insert a new row right after "Contact" (row 10) / before "Description" (row 11)
#    for the new "Jurisdiction" property, shifting the remaining rows down.
$ws.Rows.Item(11).Insert()

# Copy the formatting of the row that follows (now row 12, a normal body row)
# onto the freshly inserted blank row so it matches the rest of the table.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Fill in the new row's content.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
